# "Add files via upload" — refreshed upload of the leak-log worksheet.
# The trailing block of readings (rows 134-143) was re-entered with
# rounded figures, and the sheet's scroll position / active selection
# moved up a few rows from where the previous upload had left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 134; $r -le 143; $r++) {
    $ws.Cells.Item($r, 3).Value = 21
    $ws.Cells.Item($r, 4).Value = 133
    $ws.Cells.Item($r, 5).Value = 80
}

# Restore the view: scrolled so row 120 is at the top, D130 selected.
$win = $excel.ActiveWindow
$win.ScrollRow = 120
$win.ScrollColumn = 1
$ws.Range("D130").Select()
